# fix NPC HP error
# Insert a new "MAXHP" column before the existing "MAXMP" column (column G),
# and populate it with the same value as the HP column (F) for each NPC row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank column at G (old column G "MAXMP" and everything after
# it shifts one column to the right).
$ws.Columns.Item(7).Insert()

# Match the new column's width to its left neighbour (HP/SalePrice column),
# mirroring how the header block F:G is formatted in the saved file.
$ws.Range("F1:G1").EntireColumn.ColumnWidth = 13.29

# New header label for the inserted column.
$ws.Range("G1").Value = "MAXHP"

# Fill MAXHP with the same value as the HP (SalePrice/F) column for each data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, "A").End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $hp = $ws.Cells.Item($r, 6).Value()
    if ($hp -ne $null -and $hp -ne "") {
        $ws.Cells.Item($r, 7).Value = $hp
    }
}

# Restore the view: clear the old frozen top-left cell and move the
# selection/active cell to I8 (matches the saved workbook state).
$ws.Range("I8").Select() | Out-Null
